$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10-12: apply header label style (same as row 9) to A10:A12 ---
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10:A12").PasteSpecial(-4122) | Out-Null

# --- Score summary updates ---
$ws.Range("D10").Value = 56
$ws.Range("E10").Value = 56
$ws.Range("C11").Value = "-1"

# --- Row 15: new header cells G15/H15 (same style as existing A15/B15/D15/E15) ---
$ws.Range("E15").Copy() | Out-Null
$ws.Range("G15:H15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").Value = "Student Ans"
$ws.Range("H15").Value = "Correct Ans"

# --- Apply formats for the new G:H "question block" (rows 16-21) ---
$ws.Range("A16:A21").Copy() | Out-Null
$ws.Range("G16:G21").PasteSpecial(-4122) | Out-Null
$ws.Range("B16:B21").Copy() | Out-Null
$ws.Range("H16:H21").PasteSpecial(-4122) | Out-Null

# --- Apply formats for the extended D:E "question block" (rows 19-40) ---
$ws.Range("A19:A40").Copy() | Out-Null
$ws.Range("D19:D40").PasteSpecial(-4122) | Out-Null
$ws.Range("B19:B40").Copy() | Out-Null
$ws.Range("E19:E40").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Values for H16:H21 ---
$hvals = @{16="Option A";17="Option D";18="Option D";19="Option A";20="Option C";21="Option D"}
foreach ($r in $hvals.Keys) {
    $ws.Range("H$r").Value = $hvals[$r]
}

# --- Values for E19:E40 ---
$evals = @{
    19="Option A"; 20="Option D"; 21="Option B"; 22="Option C"; 23="Option B";
    24="Option C"; 25="Option D"; 26="Option D"; 27="Option A"; 28="Option A";
    29="Option C"; 30="Option A"; 31="Option D"; 32="Option D"; 33="Option B";
    34="Option D"; 35="Option C"; 36="Option D"; 37="Option B"; 38="Option D";
    39="Option A"; 40="Option A"
}
foreach ($r in $evals.Keys) {
    $ws.Range("E$r").Value = $evals[$r]
}
